$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for the refreshed crypto data feed.
# D-column cells are forced to Text format before assignment so that numeric-looking
# strings (e.g. "98.59", "0.0934") are preserved exactly instead of being converted
# into floating point numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.756.87"
$ws.Range("E2").Value = "  -0.03%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.294.99"
$ws.Range("E3").Value = "  -1.35%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "98.59"
$ws.Range("E5").Value = "  +3.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "270.81"
$ws.Range("E6").Value = "  -0.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.619"
$ws.Range("E7").Value = "  -1.23%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.606"
$ws.Range("E9").Value = "  -2.74%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.11"
$ws.Range("E10").Value = "  -0.66%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0934"
$ws.Range("E11").Value = "  -1.23%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.90"
$ws.Range("E12").Value = "  -2.62%  "

$ws.Range("E13").Value = "  +1.60%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.81"
$ws.Range("E14").Value = "  +1.03%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.640.76"
$ws.Range("E15").Value = "  -1.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.855"
$ws.Range("E16").Value = "  -1.07%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.295.95"
$ws.Range("E17").Value = "  -1.33%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.782.75"
$ws.Range("E18").Value = "  +0.15%  "

$ws.Range("E19").Value = "  +1.53%  "

$ws.Range("E20").Value = "  -3.13%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.22"
$ws.Range("E21").Value = "  -0.56%  "

$ws.Range("E22").Value = "  +7.86%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.72"
$ws.Range("E23").Value = "  -3.20%  "

$ws.Range("E24").Value = "  +12.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.09"
$ws.Range("E25").Value = "  -2.90%  "

$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.27"
$ws.Range("E27").Value = "  -1.42%  "

$ws.Range("E28").Value = "  -0.90%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.24"
$ws.Range("E29").Value = "  -1.55%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.82"
$ws.Range("E32").Value = "  -3.31%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0890"
$ws.Range("E33").Value = "  -1.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.42"
$ws.Range("E34").Value = "  -1.35%  "

$ws.Range("E35").Value = "  -0.22%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.73"
$ws.Range("E36").Value = "  +7.83%  "

$ws.Range("E37").Value = "  -0.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0351"
$ws.Range("E38").Value = "  -3.17%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.52"
$ws.Range("E39").Value = "  +3.97%  "

$ws.Range("E40").Value = "  -0.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.33"
$ws.Range("E41").Value = "  -1.93%  "

$ws.Range("E42").Value = "  -0.67%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.17"
$ws.Range("E43").Value = "  +0.16%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.73"
$ws.Range("E44").Value = "  +3.82%  "

$ws.Range("E45").Value = "  -3.92%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.21"
$ws.Range("E46").Value = "  -2.82%  "

$ws.Range("E47").Value = "  -1.47%  "

$ws.Range("E48").Value = "  +0.85%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "98.45"
$ws.Range("E49").Value = "  -2.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.442"
$ws.Range("E50").Value = "  +5.91%  "

$ws.Range("E51").Value = "  +11.67%  "

# Rows 30 and 31 swapped (InjectiveProtocol <-> Monero) with updated values
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "176.65"
$ws.Range("E30").Value = "  +2.14%  "

$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "38.08"
$ws.Range("E31").Value = "  -0.63%  "

